{"js": "const replacements = [\n  [\"31\u00d743=1333\", \"59\u00d748=2832\"],\n  [\"30\u00d742=1260\", \"15\u00d728=420\"],\n  [\"40\u00d778=3120\", \"86\u00d724=2064\"],\n  [\"58\u00d740=2320\", \"52\u00d726=1352\"],\n  [\"83\u00d725=2075\", \"13\u00d737=481\"],\n  [\"29\u00d727=783\", \"89\u00d765=5785\"],\n  [\"90\u00d742=3780\", \"40\u00d734=1360\"],\n  [\"14\u00d733=462\", \"29\u00d766=1914\"],\n  [\"48\u00d786=4128\", \"41\u00d714=574\"],\n  [\"50\u00d798=4900\", \"64\u00d738=2432\"],\n  [\"51\u00d715=765\", \"70\u00d721=1470\"],\n  [\"36\u00d779=2844\", \"46\u00d754=2484\"],\n  [\"70\u00d772=5040\", \"47\u00d788=4136\"],\n  [\"74\u00d711=814\", \"67\u00d719=1273\"],\n  [\"44\u00d725=1100\", \"24\u00d728=672\"],\n  [\"52\u00d774=3848\", \"74\u00d739=2886\"],\n  [\"38\u00d777=2926\", \"94\u00d739=3666\"],\n  [\"79\u00d790=7110\", \"18\u00d733=594\"],\n  [\"11\u00d740=440\", \"60\u00d779=4740\"],\n  [\"38\u00d786=3268\", \"59\u00d715=885\"],\n  [\"19\u00d753=1007\", \"75\u00d760=4500\"],\n  [\"13\u00d728=364\", \"32\u00d751=1632\"],\n  [\"27\u00d794=2538\", \"49\u00d741=2009\"],\n  [\"74\u00d784=6216\", \"87\u00d737=3219\"],\n  [\"37\u00d771=2627\", \"69\u00d747=3243\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"31\u00d743=1333\", \"59\u00d748=2832\"),\n    @(\"30\u00d742=1260\", \"15\u00d728=420\"),\n    @(\"40\u00d778=3120\", \"86\u00d724=2064\"),\n    @(\"58\u00d740=2320\", \"52\u00d726=1352\"),\n    @(\"83\u00d725=2075\", \"13\u00d737=481\"),\n    @(\"29\u00d727=783\", \"89\u00d765=5785\"),\n    @(\"90\u00d742=3780\", \"40\u00d734=1360\"),\n    @(\"14\u00d733=462\", \"29\u00d766=1914\"),\n    @(\"48\u00d786=4128\", \"41\u00d714=574\"),\n    @(\"50\u00d798=4900\", \"64\u00d738=2432\"),\n    @(\"51\u00d715=765\", \"70\u00d721=1470\"),\n    @(\"36\u00d779=2844\", \"46\u00d754=2484\"),\n    @(\"70\u00d772=5040\", \"47\u00d788=4136\"),\n    @(\"74\u00d711=814\", \"67\u00d719=1273\"),\n    @(\"44\u00d725=1100\", \"24\u00d728=672\"),\n    @(\"52\u00d774=3848\", \"74\u00d739=2886\"),\n    @(\"38\u00d777=2926\", \"94\u00d739=3666\"),\n    @(\"79\u00d790=7110\", \"18\u00d733=594\"),\n    @(\"11\u00d740=440\", \"60\u00d779=4740\"),\n    @(\"38\u00d786=3268\", \"59\u00d715=885\"),\n    @(\"19\u00d753=1007\", \"75\u00d760=4500\"),\n    @(\"13\u00d728=364\", \"32\u00d751=1632\"),\n    @(\"27\u00d794=2538\", \"49\u00d741=2009\"),\n    @(\"74\u00d784=6216\", \"87\u00d737=3219\"),\n    @(\"37\u00d771=2627\", \"69\u00d747=3243\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute(\n        [ref]$old, [ref]$true, [ref]$false, [ref]$false, [ref]$false,\n        [ref]$false, [ref]$true, 0, [ref]$false, [ref]$new, 2\n    ) | Out-Null\n}\n"}
